$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D, shifting phase_0..phase_10 (D:N) to (E:O)
$ws.Columns("D").EntireColumn.Insert()

# New column D should keep the same width/style as the column to its left (C),
# matching how Excel extends the formatted column group on insert.
$ws.Columns("D").ColumnWidth = $ws.Columns("C").ColumnWidth
$ws.Range("D1:D3").HorizontalAlignment = $ws.Range("C1").HorizontalAlignment

# New column header + values ("offset" configuration column)
$ws.Range("D1").Value = "offset"
$ws.Range("D2").Value = 0
$ws.Range("D3").Value = 0

# Restore selection to match the post-edit workbook state
$ws.Range("E2").Select()
